# #3473 replaced two properties that had gaps
# Updates the "BPS Data" sheet: replaces two buildings with gap data
# (Medstar POB North Tower -> Medstar POB South Tower; and the
# DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis
# Stevens; plus assorted owner/address corrections), fixes a trailing
# space in a header label, and clears the stray date-number-format
# that had been applied to the "Year Built" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: drop trailing space from the GHG emissions header ---
$ws.Range("M1").Value = "Total GHG Emissions Intensity"

# --- Row 2: Medstar POB North Tower -> Medstar POB South Tower ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Style = "Normal"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319
$ws.Range("M2").Value = 12.1
$ws.Range("N2").Value = 140.6

# --- Row 3: property name correction ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"
$ws.Range("I3").Style = "Normal"
$ws.Range("I3").Value = 1991

# --- Row 4: address + owner correction ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655
$ws.Range("I4").Style = "Normal"
$ws.Range("I4").Value = 1991

# --- Row 5: gross area correction ---
$ws.Range("L5").Value = 58717
$ws.Range("I5").Style = "Normal"
$ws.Range("I5").Value = 1962

# --- Row 6: President Madison Apartments -> Hampton House ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Style = "Normal"
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580
$ws.Range("M6").Value = 3.5
$ws.Range("N6").Value = 58

# --- Row 7: postal code correction ---
$ws.Range("H7").Value = 20005
$ws.Range("I7").Style = "Normal"
$ws.Range("I7").Value = 2004
$ws.Range("L7").Value = 145697

# --- Row 8: address correction ---
$ws.Range("E8").Value = "1428 H ST NW"
$ws.Range("I8").Style = "Normal"
$ws.Range("I8").Value = 1912

# --- Row 9: Year Built formatting only ---
$ws.Range("I9").Style = "Normal"
$ws.Range("I9").Value = 1880

# --- Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Style = "Normal"
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991
$ws.Range("M10").Value = 4.6
$ws.Range("N10").Value = 70.4
$ws.Range("P10").Value = 69

# --- Cosmetic: move the active selection like the source file ---
$ws.Range("D30").Select()
